$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete specific rows (by account number) that were removed from the export.
# Delete from bottom to top so row indices of earlier rows stay valid.
$rowsToDelete = @(31, 23, 22, 20, 19)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
